$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "[R] XIA(eXplainable AI) 패키지 중 DALEX로 변수 중요도 뽑기(classification)"
$ws.Range("D14").Value = "[450mm급 아두이노 드론 만들기] : 2 - 조립"
$ws.Range("D19").Value = "아기 있는 집에는 에몬스홈 그란데 가죽소파 4인"
$ws.Range("D20").Value = "프로젝트 포트폴리오 사이트 만들기 ( w/ Streamlit)"
$ws.Range("D23").Value = "온라인 얼굴고해상도 사이트 online face super-resolution api site"
$ws.Range("D24").Value = "2021년 마지막 글"
$ws.Range("E24").Value = "https://blog.naver.com/hist0134/222609768205"
$ws.Range("D26").Value = "2021 인공지능 경진대회 참가기"
$ws.Range("D28").Value = "[임피던스 제어] Improving Low-Impedance Performance"
$ws.Range("D32").Value = "Bigquery array, unnest를 mysql에서는 recursive문을 활용"
$ws.Range("E32").Value = "https://dodonam.tistory.com/345"
$ws.Range("D39").Value = "Visualize your data with Facets"
$ws.Range("D42").Value = "[BCGControlBar(MFC)] Visual Studio 2019  BCG 설치 및 개발환경 구축"
$ws.Range("D43").Value = "np.random.shuffle 과 np.random.permutation 정리"
$ws.Range("D45").Value = "Time-series forecasting at UBER"

$wb.Save()
